$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the header: LBNDIND -> LBNRIND
$ws.Range("F1").Value = "LBNRIND"

# Update the current selection/active cell as saved in the file
$ws.Range("J9").Select()
